$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: PRPC013DAAN-RC connector header (J1, J2)
$ws.Range("A13").Value = "CONN HEADER VERT 26POS 2.54MM"
$ws.Range("B13").Value = "Connector Header Through Hole 26 position 0.100"" (2.54mm)"
$ws.Range("C13").Value = "J1, J2"
$ws.Range("D13").Value = "Sullins Connector Solutions"
$ws.Range("E13").Value = "PRPC013DAAN-RC"
$ws.Range("F13").Value = "Digi-Key"
$ws.Range("H13").Value = "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PRPC013DAAN-RC/2775281"
$ws.Range("I13").Value = 1

# Row 14: S1011EC-10-ND connector header (J3)
$ws.Range("A14").Value = "CONN HEADER VERT 10POS 2.54MM"
$ws.Range("B14").Value = "Connector Header Through Hole 10 position 0.100"" (2.54mm)"
$ws.Range("C14").Value = "J3"
$ws.Range("D14").Value = "Sullins Connector Solutions"
$ws.Range("E14").Value = "S1011EC-10-ND"
$ws.Range("F14").Value = "Digi-Key"
$ws.Range("H14").Value = "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PRPC010SAAN-RC/2775244"
$ws.Range("I14").Value = 1

[void]$ws.Range("A17").Select()
